$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 229.75
$ws.Range("I12").Value = 293
$ws.Range("J12").Value = 124.333336
$ws.Range("K12").Value = 293
$ws.Range("L12").Value = 124.333336
$ws.Range("M12").Value = -123
$ws.Range("N12").Value = -464.333336
$ws.Range("H15").Value = 813.9677
$ws.Range("I15").Value = 813.9677
$ws.Range("K15").Value = 2441.9031
$ws.Range("M15").Value = -2272.9031
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H53").Value = 577.125
$ws.Range("I53").Value = 876.6
$ws.Range("J53").Value = 78
$ws.Range("K53").Value = 876.6
$ws.Range("L53").Value = 78
$ws.Range("M53").Value = -239.6
$ws.Range("N53").Value = -1352
$ws.Range("H58").Value = 1738.125
$ws.Range("J58").Value = 2065
$ws.Range("L58").Value = 6195
$ws.Range("N58").Value = -6495
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H125").Value = 8983.223
$ws.Range("J125").Value = 11249.75
$ws.Range("L125").Value = 101247.75
$ws.Range("N125").Value = -106167.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 278.30768
$ws.Range("I5").Value = 250.5
$ws.Range("J5").Value = 322.8
$ws.Range("K5").Value = 250.5
$ws.Range("L5").Value = 322.8
$ws.Range("M5").Value = -138.5
$ws.Range("N5").Value = -546.8
$ws.Range("H25").Value = 416.33334
$ws.Range("J25").Value = 200
$ws.Range("L25").Value = 200
$ws.Range("N25").Value = -1004
$ws.Range("H32").Value = 6457.85
$ws.Range("I32").Value = 5341.385
$ws.Range("K32").Value = 5341.385
$ws.Range("M32").Value = -5054.385
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976
$ws.Range("H112").Value = 35000
$ws.Range("J112").Value = 35000
$ws.Range("L112").Value = 35000
$ws.Range("N112").Value = -37954
$ws.Range("H122").Value = 2298.6
$ws.Range("I122").Value = 2298.6
$ws.Range("K122").Value = 6895.799999999999
$ws.Range("M122").Value = -4445.799999999999
$ws.Range("H132").Value = 868.5
$ws.Range("I132").Value = 868.5
$ws.Range("K132").Value = 2605.5
$ws.Range("M132").Value = -75.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 278.30768
$ws.Range("I4").Value = 250.5
$ws.Range("J4").Value = 322.8
$ws.Range("K4").Value = 250.5
$ws.Range("L4").Value = 322.8
$ws.Range("M4").Value = -135.5
$ws.Range("N4").Value = -552.8
$ws.Range("H99").Value = 1928.2858
$ws.Range("I99").Value = 1928.2858
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1928.2858
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -430.2858000000001
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 3166.6667
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3759.4
$ws.Range("I16").Value = 1712.5385
$ws.Range("J16").Value = 7560.7144
$ws.Range("K16").Value = 1712.5385
$ws.Range("L16").Value = 7560.7144
$ws.Range("M16").Value = -1425.5385
$ws.Range("N16").Value = -8134.7144
$ws.Range("H28").Value = 20630.75
$ws.Range("J28").Value = 20630.75
$ws.Range("L28").Value = 20630.75
$ws.Range("N28").Value = -21120.75
$ws.Range("H58").Value = 3632.4546
$ws.Range("I58").Value = 3196.6
$ws.Range("K58").Value = 3196.6
$ws.Range("M58").Value = -2993.6
$ws.Range("H62").Value = 7250
$ws.Range("I62").Value = 9500
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 9500
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -8876
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 7250
$ws.Range("I65").Value = 9500
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 47500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -44380
$ws.Range("N65").Value = -31240
$ws.Range("H68").Value = 37804.11
$ws.Range("J68").Value = 39996.125
$ws.Range("L68").Value = 39996.125
$ws.Range("N68").Value = -41494.125
$ws.Range("H71").Value = 37804.11
$ws.Range("J71").Value = 39996.125
$ws.Range("L71").Value = 119988.375
$ws.Range("N71").Value = -127476.375
$ws.Range("H74").Value = 38195.582
$ws.Range("J74").Value = 38195.582
$ws.Range("L74").Value = 38195.582
$ws.Range("N74").Value = -39943.582
$ws.Range("H77").Value = 38195.582
$ws.Range("J77").Value = 38195.582
$ws.Range("L77").Value = 114586.746
$ws.Range("N77").Value = -123322.746
$ws.Range("H113").Value = 3759.4
$ws.Range("I113").Value = 1712.5385
$ws.Range("J113").Value = 7560.7144
$ws.Range("K113").Value = 1712.5385
$ws.Range("L113").Value = 7560.7144
$ws.Range("M113").Value = 457.4614999999999
$ws.Range("N113").Value = -11900.7144
$ws.Range("H136").Value = 3632.4546
$ws.Range("I136").Value = 3196.6
$ws.Range("K136").Value = 9589.799999999999
$ws.Range("M136").Value = -7039.799999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 86.09090999999999
$ws.Range("J17").Value = 90
$ws.Range("L17").Value = 270
$ws.Range("N17").Value = -608
$ws.Range("H129").Value = 1567.3
$ws.Range("I129").Value = 744.1667
$ws.Range("J129").Value = 2802
$ws.Range("K129").Value = 2232.5001
$ws.Range("L129").Value = 8406
$ws.Range("M129").Value = 2767.4999
$ws.Range("N129").Value = -18406
$ws.Range("H131").Value = 2141.7058
$ws.Range("J131").Value = 2906.6667
$ws.Range("L131").Value = 8720.000100000001
$ws.Range("N131").Value = -18800.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3700
$ws.Range("I80").Value = 2450
$ws.Range("J80").Value = 3950
$ws.Range("K80").Value = 2450
$ws.Range("L80").Value = 3950
$ws.Range("M80").Value = -1452
$ws.Range("N80").Value = -5946
$ws.Range("H83").Value = 3700
$ws.Range("I83").Value = 2450
$ws.Range("J83").Value = 3950
$ws.Range("K83").Value = 12250
$ws.Range("L83").Value = 19750
$ws.Range("M83").Value = -7258
$ws.Range("N83").Value = -29734
$ws.Range("H102").Value = 1582.8
$ws.Range("I102").Value = 1228.5
$ws.Range("K102").Value = 1228.5
$ws.Range("M102").Value = 393.5
$ws.Range("H114").Value = 98333.336
$ws.Range("J114").Value = 98333.336
$ws.Range("L114").Value = 98333.336
$ws.Range("N114").Value = -107011.336
$ws.Range("H122").Value = 20875502
$ws.Range("I122").Value = 25030402
$ws.Range("K122").Value = 75091206
$ws.Range("M122").Value = -75088756
$ws.Range("H132").Value = 2997.5
$ws.Range("I132").Value = 2997.5
$ws.Range("K132").Value = 8992.5
$ws.Range("M132").Value = -6462.5
$ws.Range("H133").Value = 87997.5
$ws.Range("J133").Value = 87997.5
$ws.Range("L133").Value = 87997.5
$ws.Range("N133").Value = -98117.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1850
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1850
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1850
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2074
$ws.Range("H22").Value = 985.7143
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 985.7143
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -1714
$ws.Range("H46").Value = 3999.375
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H93").Value = 114.5
$ws.Range("I93").Value = 114.5
$ws.Range("K93").Value = 114.5
$ws.Range("M93").Value = 1133.5
$ws.Range("H110").Value = 44999.5
$ws.Range("J110").Value = 44999.5
$ws.Range("L110").Value = 44999.5
$ws.Range("N110").Value = -53179.5
$ws.Range("H132").Value = 16434.875
$ws.Range("I132").Value = 19000.6
$ws.Range("J132").Value = 12158.667
$ws.Range("K132").Value = 57001.8
$ws.Range("L132").Value = 36476.001
$ws.Range("M132").Value = -54471.8
$ws.Range("N132").Value = -41536.001
$ws.Range("H136").Value = 3367
$ws.Range("I136").Value = 3253.2
$ws.Range("J136").Value = 4505
$ws.Range("K136").Value = 9759.599999999999
$ws.Range("L136").Value = 13515
$ws.Range("M136").Value = -7209.599999999999
$ws.Range("N136").Value = -18615

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 19997.5
$ws.Range("J80").Value = 19997.5
$ws.Range("L80").Value = 19997.5
$ws.Range("N80").Value = -21993.5
$ws.Range("H83").Value = 19997.5
$ws.Range("J83").Value = 19997.5
$ws.Range("L83").Value = 59992.5
$ws.Range("N83").Value = -69976.5
$ws.Range("H132").Value = 1997
$ws.Range("I132").Value = 1997
$ws.Range("K132").Value = 5991
$ws.Range("M132").Value = -3461
$ws.Range("H136").Value = 2720.5833
$ws.Range("I136").Value = 2720.5833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8161.7499
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5611.7499
$ws.Range("N136").ClearContents()
